# runtime and sensor correction etc.
# Applies corrections to the AddrTag table:
#  - Column A: rename a few state tag names (HOTBYPASS split into 1/2, and the
#    resulting shift of RELL/RELL_HEAT/NORMAL_WATER/COLD_WATER/RSV_BPOS_1/RSV_BPOS_2)
#  - Column J: prefix DO remark text with "DO " and renumber the "预留" placeholders
#  - Move active selection to M29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column A corrections (state tag names) ----
$ws.Range("A15").Value = "state_DO_HOTBYPASS1_BPOS"
$ws.Range("A16").Value = "state_DO_HOTBYPASS2_BPOS"
$ws.Range("A17").Value = "state_DO_RELL_BPOS"
$ws.Range("A18").Value = "state_DO_RELL_HEAT_BPOS"
$ws.Range("A19").Value = "state_DO_NORMAL_WATER_BPOS"
$ws.Range("A20").Value = "state_DO_COLD_WATER_BPOS"
$ws.Range("A21").Value = "state_DO_RSV_BPOS_1"
$ws.Range("A29").Value = "state_DO_RSV_BPOS_2"

# ---- Column J corrections (remark text) ----
$ws.Range("J3").Value  = "DO 沉淀滤芯电磁阀"
$ws.Range("J4").Value  = "DO 活性炭滤芯电磁阀"
$ws.Range("J5").Value  = "DO 饮水箱进水"
$ws.Range("J6").Value  = "DO 循环支路"
$ws.Range("J7").Value  = "DO 沉淀杯清洗"
$ws.Range("J8").Value  = "DO 膨胀水箱"
$ws.Range("J9").Value  = "DO 冰水阀"
$ws.Range("J10").Value = "DO 循环泵"
$ws.Range("J11").Value = "DO 气泵"
$ws.Range("J12").Value = "DO UV"
$ws.Range("J13").Value = "DO 压机1启动"
$ws.Range("J14").Value = "DO 压机2"
$ws.Range("J15").Value = "DO 热气旁通"
$ws.Range("J16").Value = "DO 热气旁通2"
$ws.Range("J17").Value = "DO 转轮"
$ws.Range("J18").Value = "DO 转轮电加热"
$ws.Range("J19").Value = "DO 常温水"
$ws.Range("J20").Value = "DO 冰水"
$ws.Range("J21").Value = "DO 预留1"
$ws.Range("J22").Value = "DO 预留"
$ws.Range("J23").Value = "DO 预留"
$ws.Range("J24").Value = "DO 预留"
$ws.Range("J25").Value = "DO 预留"
$ws.Range("J26").Value = "DO 绿灯"
$ws.Range("J27").Value = "DO 红灯"
$ws.Range("J28").Value = "DO 橙灯"
$ws.Range("J29").Value = "DO 预留2"
$ws.Range("J30").Value = "DO 预留"
$ws.Range("J31").Value = "DO 预留5"
$ws.Range("J32").Value = "DO 预留6"
$ws.Range("J33").Value = "DO 预留7"
$ws.Range("J34").Value = "DO 预留8"

# ---- Update the active selection shown in the sheet view ----
$ws.Range("M29").Select()
